# Refresh the crypto price/volume table (columns D and E) with the latest
# scraped values, matching the GitHub Actions update commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '43.424.71'
$ws.Range('E2').Value = '  +2.86%  '
$ws.Range('D3').Value = '2.309.08'
$ws.Range('E3').Value = '  +1.71%  '
# Leading apostrophe forces text so Excel doesn't drop the trailing zero
# and turn "1.00" into the number 1.
$ws.Range('D4').Value = '''1.00'
$ws.Range('D5').Value = '310.98'
$ws.Range('E5').Value = '  +1.40%  '
$ws.Range('D6').Value = '103.22'
$ws.Range('E6').Value = '  +5.82%  '
$ws.Range('E7').Value = '  +1.61%  '
$ws.Range('E8').Value = '  -0.02%  '
$ws.Range('D9').Value = '0.531'
$ws.Range('E9').Value = '  +8.28%  '
$ws.Range('D10').Value = '35.75'
$ws.Range('E10').Value = '  +1.24%  '
$ws.Range('D11').Value = '0.0813'
$ws.Range('E11').Value = '  +3.06%  '
$ws.Range('E12').Value = '  -0.96%  '
$ws.Range('D13').Value = '7.03'
$ws.Range('E13').Value = '  +2.43%  '
$ws.Range('D14').Value = '2.665.91'
$ws.Range('E14').Value = '  +1.71%  '
$ws.Range('D15').Value = '15.03'
$ws.Range('E15').Value = '  +2.10%  '
$ws.Range('D16').Value = '2.413.39'
$ws.Range('E16').Value = '  +6.56%  '
$ws.Range('D17').Value = '0.809'
$ws.Range('E17').Value = '  +2.25%  '
$ws.Range('D18').Value = '43.324.38'
$ws.Range('E18').Value = '  +2.90%  '
$ws.Range('D19').Value = '12.28'
$ws.Range('E19').Value = '  -0.06%  '
$ws.Range('E20').Value = '  +3.28%  '
$ws.Range('D21').Value = '6.18'
$ws.Range('E21').Value = '  +2.92%  '
$ws.Range('D22').Value = '68.12'
$ws.Range('E22').Value = '  +0.59%  '
$ws.Range('D23').Value = '241.55'
$ws.Range('E23').Value = '  +1.85%  '
$ws.Range('D24').Value = '2.02'
$ws.Range('E24').Value = '  +1.49%  '
$ws.Range('D25').Value = '2.62'
$ws.Range('E25').Value = '  +1.35%  '
$ws.Range('E26').Value = '  +0.05%  '
$ws.Range('D27').Value = '24.99'
$ws.Range('E27').Value = '  +6.21%  '
$ws.Range('E28').Value = '  +8.02%  '
$ws.Range('D29').Value = '36.84'
$ws.Range('E29').Value = '  -1.40%  '
$ws.Range('D30').Value = '9.67'
$ws.Range('E30').Value = '  +0.82%  '
$ws.Range('D31').Value = '171.78'
$ws.Range('E31').Value = '  +5.49%  '
$ws.Range('E32').Value = '  +0.40%  '
$ws.Range('D33').Value = '0.999'
$ws.Range('E33').Value = '  -0.08%  '
$ws.Range('D34').Value = '2.55'
$ws.Range('E34').Value = '  +6.88%  '
$ws.Range('D35').Value = '17.85'
$ws.Range('E35').Value = '  +0.82%  '
$ws.Range('D36').Value = '0.0742'
$ws.Range('E36').Value = '  +0.98%  '
$ws.Range('D37').Value = '3.07'
$ws.Range('E37').Value = '  -2.15%  '
# Rows 38/39 swap places: Kaspa moves up to row 38, ARBITRUM moves down
# to row 39 (ranking shuffled), each bringing its own updated price/volume.
$ws.Range('B38').Value = 'Kaspa'
$ws.Range('C38').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D38').Value = '0.107'
$ws.Range('E38').Value = '  +2.85%  '
$ws.Range('B39').Value = 'ARBITRUM'
$ws.Range('C39').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D39').Value = '1.89'
$ws.Range('E39').Value = '  +3.30%  '
$ws.Range('E40').Value = '  +1.68%  '
$ws.Range('E41').Value = '  +5.31%  '
$ws.Range('D42').Value = '2.31'
$ws.Range('E42').Value = '  -1.08%  '
$ws.Range('E43').Value = '  +4.37%  '
$ws.Range('D44').Value = '1.973.66'
$ws.Range('E44').Value = '  +1.27%  '
$ws.Range('D45').Value = '19.19'
$ws.Range('E45').Value = '  +0.83%  '
# Same trailing-zero-preservation trick as D4 above.
$ws.Range('D46').Value = '''3.00'
$ws.Range('E46').Value = '  +2.96%  '
$ws.Range('D47').Value = '9.95'
$ws.Range('E47').Value = '  -0.23%  '
$ws.Range('D48').Value = '55.56'
$ws.Range('E48').Value = '  +3.39%  '
$ws.Range('E49').Value = '  +1.85%  '
$ws.Range('D50').Value = '1.59'
$ws.Range('E50').Value = '  +7.78%  '
$ws.Range('D51').Value = '2.532.32'
$ws.Range('E51').Value = '  +1.57%  '
